# Apply cryptos list price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.682.79'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '2.199.35'
$ws.Range('E3').Value = '  -2.91%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'229.11"
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('D6').Value = "'0.615"
$ws.Range('E6').Value = '  -4.31%  '
$ws.Range('D7').Value = "'60.00"
$ws.Range('E7').Value = '  -5.00%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = "'0.399"
$ws.Range('E9').Value = '  -2.76%  '
$ws.Range('D10').Value = "'56.88"
$ws.Range('E10').Value = '  -5.09%  '
$ws.Range('D11').Value = "'0.0882"
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('E12').Value = '  -2.20%  '
$ws.Range('D13').Value = '2.527.22'
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('D14').Value = "'15.26"
$ws.Range('E14').Value = '  -4.81%  '
$ws.Range('D15').Value = "'22.05"
$ws.Range('E15').Value = '  -3.40%  '
$ws.Range('D16').Value = "'0.790"
$ws.Range('E16').Value = '  -3.77%  '
$ws.Range('D17').Value = "'5.54"
$ws.Range('E17').Value = '  -3.13%  '
$ws.Range('D18').Value = '2.195.56'
$ws.Range('E18').Value = '  -3.23%  '
$ws.Range('D19').Value = '41.567.96'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').Value = "'71.80"
$ws.Range('E20').Value = '  -3.85%  '
$ws.Range('D21').Value = '0.0₃0896'
$ws.Range('E21').Value = '  -3.36%  '
$ws.Range('D22').Value = "'6.00"
$ws.Range('E22').Value = '  -2.44%  '
$ws.Range('D23').Value = "'240.66"
$ws.Range('E23').Value = '  -5.22%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -3.62%  '
$ws.Range('E26').Value = '  -3.39%  '
$ws.Range('D27').Value = "'9.55"
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('D28').Value = "'167.78"
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('D29').Value = "'0.139"
$ws.Range('E29').Value = '  -6.44%  '
$ws.Range('D30').Value = "'1.44"
$ws.Range('E30').Value = '  -0.97%  '
$ws.Range('D31').Value = "'19.65"
$ws.Range('E31').Value = '  -4.34%  '
$ws.Range('D32').Value = "'2.58"
$ws.Range('E32').Value = '  -8.95%  '
$ws.Range('E33').Value = '  -3.66%  '
$ws.Range('D34').Value = "'4.97"
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').Value = "'4.56"
$ws.Range('E35').Value = '  -4.93%  '
$ws.Range('D36').Value = "'0.0640"
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').Value = "'3.54"
$ws.Range('E37').Value = '  -8.23%  '
$ws.Range('D38').Value = "'6.25"
$ws.Range('E38').Value = '  -7.89%  '
$ws.Range('D39').Value = "'2.32"
$ws.Range('E39').Value = '  -5.44%  '
$ws.Range('B40').Value = 'TerraClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D40').Value = "'0.000242"
$ws.Range('E40').Value = '  -1.66%  '
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D41').Value = "'0.999"
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('D42').Value = "'0.0239"
$ws.Range('D43').Value = "'8.60"
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = "'0.0952"
$ws.Range('E44').Value = '  -2.87%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = "'4.39"
$ws.Range('E45').Value = '  -13.38%  '
$ws.Range('E46').Value = '  -3.66%  '
$ws.Range('D47').Value = "'95.89"
$ws.Range('E47').Value = '  -5.90%  '
$ws.Range('D48').Value = '1.453.63'
$ws.Range('E48').Value = '  -3.31%  '
$ws.Range('E49').Value = '  -1.84%  '
$ws.Range('D50').Value = "'16.04"
$ws.Range('E50').Value = '  -10.11%  '
$ws.Range('B51').Value = 'HuobiBTC'
$ws.Range('C51').Value = 'https://coinranking.com/coin/upmyKdAQ+huobibtc-hbtc'
$ws.Range('D51').Value = '151.218.32'
$ws.Range('E51').Value = '  +301.16%  '
